$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+1+Marron.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+10+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+11+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+12+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+2+Marron.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+3+Marron.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+4+Marron.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+5+Morado.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+6+Morado.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+7+Morado.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+8+Morado.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Camiseta+More+than/Camiseta+_More+than_+9+Rosa.jpg']"
$ws.Range("O4").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Top+Toxic+Love/Top+Toxic+Love+1+Blanco.jpg']"
$ws.Range("O5").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Top+The+Funky+Love/Top+The+Funky+Love+1+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Top+The+Funky+Love/Top+The+Funky+Love+2+Estampado.jpg']"
$ws.Range("O6").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Pantal%E0%B8%82n+The+Funky+Love/Pantalon+The+Funky+Love+1+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Pantal%E0%B8%82n+The+Funky+Love/Pantalon+The+Funky+Love+2+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Pantal%E0%B8%82n+The+Funky+Love/Pantalon+The+Funky+Love+3+Estampado.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Pantal%E0%B8%82n+The+Funky+Love/Pantalon+The+Funky+Love+4+Estampado.jpg']"
$ws.Range("O7").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Pantal%E0%B8%82n+The+Jungle+love/Pantalon+The+Jungle+love+1+Estampado.jpg']"
$ws.Range("O8").Value = "['https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Falda+The+Psychedelic+Love/Falda+The+Psychedelic+Love++3+Morado.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Falda+The+Psychedelic+Love/Falda+The+Psychedelic+Love++5+Rosa.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Falda+The+Psychedelic+Love/Falda+The+Psychedelic+Love+1+Azul.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Falda+The+Psychedelic+Love/Falda+The+Psychedelic+Love+2+Azul.jpg', 'https://recursosmolova.s3.amazonaws.com/Products+Images/Madly+In+Love/Falda+The+Psychedelic+Love/Falda+The+Psychedelic+Love+4+Morado.jpg']"
